$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.525.70"
$ws.Range("E2").Value = "  +1.67%  "

$ws.Range("D3").Value = "2.478.29"
$ws.Range("E3").Value = "  +1.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.64"
$ws.Range("E5").Value = "  +1.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.28"
$ws.Range("E6").Value = "  +1.32%  "

$ws.Range("E7").Value = "  -0.84%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +3.15%  "

$ws.Range("E10").Value = "  -0.42%  "

$ws.Range("E11").Value = "  +1.99%  "

$ws.Range("E12").Value = "  +3.20%  "

$ws.Range("D13").Value = "2.862.69"
$ws.Range("E13").Value = "  +1.36%  "

$ws.Range("E14").Value = "  -0.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.04"
$ws.Range("E15").Value = "  +9.64%  "

$ws.Range("D16").Value = "2.471.06"
$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("E17").Value = "  -0.99%  "

$ws.Range("D18").Value = "41.522.66"
$ws.Range("E18").Value = "  +1.67%  "

$ws.Range("E19").Value = "  +2.86%  "

$ws.Range("D20").Value = "0.0₃0936"
$ws.Range("E20").Value = "  +3.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.70"
$ws.Range("E21").Value = "  +6.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.30"
$ws.Range("E22").Value = "  +3.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.57"
$ws.Range("E23").Value = "  +1.50%  "

$ws.Range("E24").Value = "  -0.47%  "

$ws.Range("E25").Value = "  -0.42%  "

$ws.Range("E26").Value = "  +1.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.84"
$ws.Range("E27").Value = "  +5.82%  "

$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("E29").Value = "  +1.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.90"
$ws.Range("E30").Value = "  +1.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.92"
$ws.Range("E31").Value = "  +5.20%  "

$ws.Range("E32").Value = "  +1.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.57"
$ws.Range("E33").Value = "  +1.63%  "

$ws.Range("E34").Value = "  +3.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.39"
$ws.Range("E35").Value = "  +4.87%  "

$ws.Range("E36").Value = "  -7.65%  "

$ws.Range("E37").Value = "  -0.82%  "

$ws.Range("E38").Value = "  +4.52%  "

$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("E40").Value = "  +0.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.11"
$ws.Range("E41").Value = "  +0.29%  "

$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.73"
$ws.Range("E43").Value = "  +0.93%  "

$ws.Range("D44").Value = "1.971.57"
$ws.Range("E44").Value = "  +0.88%  "

$ws.Range("E45").Value = "  +1.72%  "

$ws.Range("E46").Value = "  -0.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.14"
$ws.Range("E47").Value = "  +7.68%  "

$ws.Range("D48").Value = "2.721.31"
$ws.Range("E48").Value = "  +1.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.74"
$ws.Range("E49").Value = "  +2.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.16"
$ws.Range("E50").Value = "  -0.32%  "

$ws.Range("E51").Value = "  -1.37%  "
